# "Droit de mutation dans spécification"
#
# Adds a "Droit de mutation - Taxe de Bienvenue" analysis block to the
# "montage financier" sheet (land-purchase cost table + 5 hypothesis
# scenarios H-1..H-5 with explanations), fixes a typo in "spécification",
# and updates a couple of pre-existing values.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # montage financier
$ws3 = $wb.Worksheets.Item(3)   # spécification

# ---------------------------------------------------------------------
# 1. spécification sheet: fix a typo in-place, then append a new note.
# ---------------------------------------------------------------------
$ws3.Range("A13").Value = "lumière extérieure"
$ws3.Range("A14").Value = "vérifier quand ikea spéciaux 20% cuisine au Canada"

# ---------------------------------------------------------------------
# 2. montage financier: existing hypothèque line, corrected amount, plus
#    a new interest-rate / cap column layout (D blank spacer, E rate,
#    F cap) alongside the original A/C columns.
# ---------------------------------------------------------------------
$ws1.Range("C4").Value = 97500
$ws1.Range("E4").Formula = "=5/1000"
$ws1.Range("F4").Value = 50000

$ws1.Range("C5").Value = 13000
$ws1.Range("E5").Value = 0.01
$ws1.Range("F5").Value = 250000

$ws1.Range("C6").Value = 81000
$ws1.Range("E6").Formula = "=1.5/100"
$ws1.Range("F6").Value = 500000

$ws1.Range("C7").Value = 110000

$ws1.Range("C8").Value = 10000

$ws1.Range("C12").Formula = "=SUM(C4:C8)"

# ---------------------------------------------------------------------
# 3. Land values table (written before the H-1..H-5 header so the
#    shared-string table picks up "terrain …" / "Bienvenue …" first).
# ---------------------------------------------------------------------
$ws1.Range("A19").Value = "terrain Robert"
$ws1.Range("B19").Value = 21632
$ws1.Range("D19").Formula = "=C20/B20 *B19"
$ws1.Range("E19").Value = 50000
$ws1.Range("F19").Value = 50000

$ws1.Range("A20").Value = "terrain Boutin"
$ws1.Range("B20").Value = 214400
$ws1.Range("C20").Value = 325000
$ws1.Range("D20").Formula = "=C20-D19"
$ws1.Range("E20").Value = 275000
$ws1.Range("F20").Value = 265000

# ---------------------------------------------------------------------
# 4. Bienvenue (land-transfer tax) computations.
# ---------------------------------------------------------------------
$ws1.Range("A22").Value = "Bienvenue Robert"
$ws1.Range("B22").Formula = "=B19 *`$E4"
$ws1.Range("D22").Formula = "=D19 *`$E4"
$ws1.Range("E22").Formula = "=E19 *`$E4"
$ws1.Range("F22").Formula = "=F19 *`$E4"

$ws1.Range("A24").Value = "Bienvenue totale"
$ws1.Range("B24").Formula = "=B22+B23"
$ws1.Range("C24").Formula = "=C22+C23"
$ws1.Range("D24").Formula = "=D22+D23"
$ws1.Range("E24").Formula = "=E22+E23"
$ws1.Range("F24").Formula = "=F22+F23"

$ws1.Range("A23").Value = "Bienvenue Boutin"
$ws1.Range("B23").Formula = "=(`$E4 *`$F4) + (B20 -`$F4) *`$E5"
$ws1.Range("C23").Formula = "=(`$E4 *`$F4) + (C20 -`$F4) *`$E5"
$ws1.Range("D23").Formula = "=(`$E4 *`$F4) + (D20 -`$F4) * `$E5"
$ws1.Range("E23").Formula = "=(`$E4 *`$F4) + (E20 -`$F4) * `$E5"
$ws1.Range("F23").Formula = "=(`$E4 *`$F4) + (F20 -`$F4) * `$E5"

# ---------------------------------------------------------------------
# 5. Header row (bold, centered) H-1..H-5.
# ---------------------------------------------------------------------
$ws1.Range("B18").Value = "H-1"
$ws1.Range("C18").Value = "H-2"
$ws1.Range("D18").Value = "H-3"
$ws1.Range("E18").Value = "H-4"
$ws1.Range("F18").Value = "H-5"
$ws1.Range("A18:F18").Font.Bold = $true
$ws1.Range("B18:F18").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 6. Hypothesis explanations H-1..H-5.
# ---------------------------------------------------------------------
$ws1.Range("A27").Value = "H-1"
$ws1.Range("B27").Value = "Achat au prix de l'évaluation"

$ws1.Range("A28").Value = "H-2"
$ws1.Range("B28").Value = "Achat au prix réel tout mis sur le terrain de M. Boutin"

$ws1.Range("A29").Value = "H-3"
$ws1.Range("B29").Value = "Achat au prix réel et le montant est proportionné selon l'évaluation de chaque terrain "

$ws1.Range("A30").Value = "H-4"
$ws1.Range("B30").Value = "Achat au prix réel et le prix d'achat du terrain de Mme Robert à 50 000$"

$ws1.Range("A31").Value = "H-5"
$ws1.Range("B31").Value = "Comme H-4 mais au prix total de 315 000$"

# ---------------------------------------------------------------------
# 7. New "Droit de mutation" section title (last new string -> index 40).
# ---------------------------------------------------------------------
$ws1.Range("A16").Value = "Droit de mutation - Taxe de Bienvenue"

# ---------------------------------------------------------------------
# 8. Number formats: currency on the amount grid, percentage on rates.
# ---------------------------------------------------------------------
$ws1.Range("B4:F17").NumberFormat = "_-[`$`$-409]* #,##0.00_ ;_-[`$`$-409]* \-#,##0.00\ ;_-[`$`$-409]* ""-""??_ ;_-@_ "
$ws1.Range("B18:F18").NumberFormat = "_-[`$`$-409]* #,##0.00_ ;_-[`$`$-409]* \-#,##0.00\ ;_-[`$`$-409]* ""-""??_ ;_-@_ "
$ws1.Range("B19:F24").NumberFormat = "_-[`$`$-409]* #,##0.00_ ;_-[`$`$-409]* \-#,##0.00\ ;_-[`$`$-409]* ""-""??_ ;_-@_ "
$ws1.Range("E4:E6").NumberFormat = "0.0%"

# ---------------------------------------------------------------------
# 9. Column widths.
# ---------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 20.43
$ws1.Range("B1:F1").EntireColumn.ColumnWidth = 11.42

# ---------------------------------------------------------------------
# 10. View state: "montage financier" becomes the active/selected tab,
#     "spécification" loses tabSelected and gets a fresh selection.
# ---------------------------------------------------------------------
$ws3.Range("A16").Select()
$ws1.Activate()
$ws1.Range("C28").Select()
